$wb = $excel.ActiveWorkbook

# APPL was a typo that should always have been AAPL - fix every occurrence
$ws_rsu = $wb.Worksheets.Item("rsu")
$ws_rsu.Range("B6").Value = "AAPL"

$ws_dividends = $wb.Worksheets.Item("dividends")
$ws_dividends.Range("B3").Value = "AAPL"

$ws_sell = $wb.Worksheets.Item("sell_orders")
$ws_sell.Range("B6").Value = "AAPL"
$ws_sell.Range("B7").Value = "AAPL"

# restore the active/selected sheet back to the first one (rsu)
$ws_rsu.Activate()
